$d = $word.ActiveDocument

# Locate the exact literal text "DATE:07/10/2023" (the day "07" is being
# corrected to "10"); Find.Execute returns a Range collapsed onto the hit,
# which does not include the paragraph mark, so the paragraph's own
# attributes (paraId/rsid/etc.) are left untouched by the XML splice below.
$rng = $d.Content
$found = $rng.Find.Execute("DATE:07/10/2023", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found -and $rng.Find.Found) {
    # Replace the matched run's text with three separate runs - "DATE:",
    # "10" and "/10/2023" - by splicing in a WordOpenXML fragment. This
    # reproduces exactly the run split produced when the date's day value
    # is edited in place, without stamping any incidental formatting onto
    # the newly created runs.
    $xmlSnippet = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>DATE:</w:t></w:r><w:r><w:t>10</w:t></w:r><w:r><w:t>/10/2023</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    $rng.InsertXML($xmlSnippet)
}
